$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in row 10 (8th data row of the TEAMS table): GetDirectoryStructure API ---
# Shared strings must be appended in this exact order (Response, Body, URL) to line up
# with the target workbook's sharedStrings.xml, so set the cells in that sequence.
$response10 = "[" + "`n" + "    ""index.html""," + "`n" + "    ""scripts.js""" + "`n" + "]"
$body10     = "{" + "`n" + "    ""Path"":""~/Projects/team1""" + "`n" + "}"

$ws.Range("E10").Value = $response10
$ws.Range("D10").Value = $body10
$ws.Range("C10").Value = "api/GetDirectoryStructure"
$ws.Range("B10").Value = "POST"

# Body/Response columns wrap their (multi-line) text, like the other rows in the table.
$ws.Range("D10").WrapText = $true
$ws.Range("E10").WrapText = $true

# Row grows to fit the new multi-line content.
$ws.Rows.Item(10).RowHeight = 60

# --- Move the active selection, as recorded in the saved sheet view ---
$ws.Range("E10").Select()
